# Adds a new "Intermediate SQL" course entry right after the existing
# "Introduction to SQL" row (row 12), shifting the remaining rating rows
# down by one, and updates the selected cell accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 13 (pushes old row 13.. down to 14..)
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the course name and its "disk_savvy" rating
$ws.Range("A13").Value() = "Intermediate SQL"
$ws.Range("G13").Value() = 5

# Reflect the new active selection in the sheet view
$ws.Range("G14").Select()
